# Auto-generated edit script: applies scheduled-runner market data refresh
# to the Maduin_Profits workbook, updating H-N (price/profit) columns
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1750
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 1750
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H70").Value = 1561.3334
$ws.Range("I70").Value = 780.4
$ws.Range("K70").Value = 2341.2
$ws.Range("M70").Value = -2071.2
$ws.Range("H73").Value = 1561.3334
$ws.Range("I73").Value = 780.4
$ws.Range("K73").Value = 2341.2
$ws.Range("M73").Value = -1405.2
$ws.Range("H86").Value = 6332.1665
$ws.Range("I86").Value = 4999.5
$ws.Range("K86").Value = 4999.5
$ws.Range("M86").Value = -3876.5
$ws.Range("H89").Value = 6332.1665
$ws.Range("I89").Value = 4999.5
$ws.Range("K89").Value = 24997.5
$ws.Range("M89").Value = -19381.5
$ws.Range("H111").Value = 1154.3334
$ws.Range("J111").Value = 1275.6666
$ws.Range("L111").Value = 3826.9998
$ws.Range("N111").Value = -9960.9998
$ws.Range("H115").Value = 541
$ws.Range("I115").Value = 449.2
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 1347.6
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = 219.4000000000001
$ws.Range("N115").Value = -6134
$ws.Range("H116").Value = 4723.8
$ws.Range("I116").Value = 4496.8
$ws.Range("K116").Value = 4496.8
$ws.Range("M116").Value = -1054.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2225.7273
$ws.Range("I2").Value = 1053.8889
$ws.Range("J2").Value = 7499
$ws.Range("K2").Value = 1053.8889
$ws.Range("L2").Value = 7499
$ws.Range("M2").Value = -940.8888999999999
$ws.Range("N2").Value = -7725
$ws.Range("H74").Value = 1378.0834
$ws.Range("I74").Value = 1378.0834
$ws.Range("K74").Value = 1378.0834
$ws.Range("M74").Value = -504.0834
$ws.Range("H77").Value = 1378.0834
$ws.Range("I77").Value = 1378.0834
$ws.Range("K77").Value = 6890.416999999999
$ws.Range("M77").Value = -2522.416999999999
$ws.Range("H110").Value = 424.85715
$ws.Range("I110").Value = 329
$ws.Range("K110").Value = 329
$ws.Range("M110").Value = 1716
$ws.Range("H116").Value = 2225.7273
$ws.Range("I116").Value = 1053.8889
$ws.Range("J116").Value = 7499
$ws.Range("K116").Value = 1053.8889
$ws.Range("L116").Value = 7499
$ws.Range("M116").Value = 1240.1111
$ws.Range("N116").Value = -12087
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2225.7273
$ws.Range("I3").Value = 1053.8889
$ws.Range("J3").Value = 7499
$ws.Range("K3").Value = 1053.8889
$ws.Range("L3").Value = 7499
$ws.Range("M3").Value = -939.8888999999999
$ws.Range("N3").Value = -7727
$ws.Range("H81").Value = 32500
$ws.Range("J81").Value = 32500
$ws.Range("L81").Value = 32500
$ws.Range("N81").Value = -34622
$ws.Range("H84").Value = 32500
$ws.Range("J84").Value = 32500
$ws.Range("L84").Value = 97500
$ws.Range("N84").Value = -108108
$ws.Range("H99").Value = 2653.5
$ws.Range("I99").Value = 1952.5
$ws.Range("J99").Value = 3354.5
$ws.Range("K99").Value = 1952.5
$ws.Range("L99").Value = 3354.5
$ws.Range("M99").Value = -454.5
$ws.Range("N99").Value = -6350.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1816
$ws.Range("I31").Value = 1849
$ws.Range("K31").Value = 1849
$ws.Range("M31").Value = -1554
$ws.Range("H34").Value = 1816
$ws.Range("I34").Value = 1849
$ws.Range("K34").Value = 1849
$ws.Range("M34").Value = -1647
$ws.Range("H107").Value = 2123.2727
$ws.Range("I107").Value = 1642.625
$ws.Range("J107").Value = 3405
$ws.Range("K107").Value = 1642.625
$ws.Range("L107").Value = 3405
$ws.Range("M107").Value = 277.375
$ws.Range("N107").Value = -7245

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 988.5
$ws.Range("I18").Value = 986.2
$ws.Range("K18").Value = 2958.6
$ws.Range("M18").Value = -2789.6
$ws.Range("H21").Value = 750
$ws.Range("I21").Value = 750
$ws.Range("K21").Value = 2250
$ws.Range("M21").Value = -2077
$ws.Range("H59").Value = 1000
$ws.Range("I59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -2460
$ws.Range("H97").Value = 880.5
$ws.Range("I97").Value = 765.6667
$ws.Range("K97").Value = 2297.0001
$ws.Range("M97").Value = -1801.0001
$ws.Range("H131").Value = 863.44446
$ws.Range("I131").Value = 442
$ws.Range("J131").Value = 1074.1666
$ws.Range("K131").Value = 1326
$ws.Range("L131").Value = 3222.4998
$ws.Range("M131").Value = 3714
$ws.Range("N131").Value = -13302.4998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12372
$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41856
$ws.Range("H132").Value = 3949.6
$ws.Range("I132").Value = 3856.7144
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 11570.1432
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -9040.143199999999
$ws.Range("N132").Value = -17559.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12125
$ws.Range("H27").Value = 12125
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1570
$ws.Range("J107").Value = 1350
$ws.Range("L107").Value = 4050
$ws.Range("N107").Value = -7890
$ws.Range("H113").Value = 741
$ws.Range("J113").Value = 565.625
$ws.Range("L113").Value = 1696.875
$ws.Range("N113").Value = -6036.875
$ws.Range("H136").Value = 887.53845
$ws.Range("I136").Value = 887.53845
$ws.Range("K136").Value = 2662.61535
$ws.Range("M136").Value = -112.61535

Write-Output "Applied scheduled market-data refresh to 8 sheets (169 cell updates)."